$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $escaped = $val.Replace('"', '""')
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

$ws.Range("D2").Value2 = "65.134.32"
$ws.Range("E2").Value2 = "  +2.18%  "
$ws.Range("D3").Value2 = "3.163.43"
$ws.Range("E3").Value2 = "  +3.74%  "
Set-TextValue "D5" '578.66'
$ws.Range("E5").Value2 = "  +4.07%  "
Set-TextValue "D6" '150.81'
$ws.Range("E6").Value2 = "  +6.47%  "
$ws.Range("E7").Value2 = "  +0.09%  "
$ws.Range("D8").Value2 = "3.159.35"
$ws.Range("E8").Value2 = "  +3.72%  "
Set-TextValue "D9" '0.530'
$ws.Range("E9").Value2 = "  +2.06%  "
$ws.Range("E10").Value2 = "  +5.98%  "
$ws.Range("E11").Value2 = "  -0.01%  "
Set-TextValue "D12" '0.502'
$ws.Range("E12").Value2 = "  +5.15%  "
Set-TextValue "D13" '0.0000271'
$ws.Range("E13").Value2 = "  +17.04%  "
Set-TextValue "D14" '37.47'
$ws.Range("E14").Value2 = "  +6.73%  "
$ws.Range("D15").Value2 = "3.687.03"
$ws.Range("E15").Value2 = "  +3.95%  "
$ws.Range("D16").Value2 = "65.225.17"
$ws.Range("E16").Value2 = "  +2.27%  "
$ws.Range("D17").Value2 = "3.170.82"
$ws.Range("E17").Value2 = "  +3.95%  "
Set-TextValue "D18" '7.17'
$ws.Range("E18").Value2 = "  +6.29%  "
$ws.Range("E19").Value2 = "  +1.33%  "
Set-TextValue "D20" '511.01'
$ws.Range("E20").Value2 = "  +4.89%  "
Set-TextValue "D21" '14.85'
$ws.Range("E21").Value2 = "  +5.22%  "
Set-TextValue "D22" '0.725'
$ws.Range("E22").Value2 = "  +6.46%  "
Set-TextValue "D23" '15.30'
$ws.Range("E23").Value2 = "  +6.37%  "
Set-TextValue "D24" '7.81'
$ws.Range("E24").Value2 = "  +4.07%  "
Set-TextValue "D25" '85.02'
$ws.Range("E25").Value2 = "  +3.11%  "
Set-TextValue "D26" '0.999'
$ws.Range("E26").Value2 = "  -0.08%  "
Set-TextValue "D27" '9.08'
$ws.Range("E27").Value2 = "  +12.23%  "
Set-TextValue "D28" '2.93'
$ws.Range("E28").Value2 = "  +5.01%  "
Set-TextValue "D29" '2.19'
$ws.Range("E29").Value2 = "  +8.24%  "
Set-TextValue "D30" '2.82'
$ws.Range("E30").Value2 = "  +15.34%  "
Set-TextValue "D31" '27.82'
$ws.Range("E31").Value2 = "  +6.06%  "
$ws.Range("E32").Value2 = "  +0.09%  "
$ws.Range("E33").Value2 = "  +4.06%  "
$ws.Range("E34").Value2 = "  +11.72%  "
Set-TextValue "D35" '6.60'
$ws.Range("E35").Value2 = "  +6.75%  "
Set-TextValue "D36" '55.78'
$ws.Range("E36").Value2 = "  +1.07%  "
Set-TextValue "D37" '0.0901'
$ws.Range("E37").Value2 = "  +10.85%  "
Set-TextValue "D38" '474.31'
$ws.Range("E38").Value2 = "  +7.86%  "
Set-TextValue "D39" '3.06'
$ws.Range("E39").Value2 = "  +10.99%  "
Set-TextValue "D40" '0.0420'
$ws.Range("E40").Value2 = "  +3.17%  "
Set-TextValue "D41" '8.66'
$ws.Range("E41").Value2 = "  +4.31%  "
$ws.Range("D42").Value2 = "3.072.56"
$ws.Range("E42").Value2 = "  +1.88%  "
Set-TextValue "D43" '0.118'
$ws.Range("E43").Value2 = "  +2.71%  "
Set-TextValue "D44" '2.46'
$ws.Range("E44").Value2 = "  +11.37%  "
Set-TextValue "D45" '0.285'
$ws.Range("E45").Value2 = "  +5.69%  "
Set-TextValue "D46" '29.18'
$ws.Range("E46").Value2 = "  +5.79%  "
$ws.Range("D47").Value2 = "0.0₃0603"
$ws.Range("E47").Value2 = "  +17.94%  "
$ws.Range("E49").Value2 = "  +1.31%  "
Set-TextValue "D50" '2.25'
$ws.Range("E50").Value2 = "  +8.07%  "
Set-TextValue "D51" '120.38'
$ws.Range("E51").Value2 = "  +2.22%  "
